# Delete unused time table
# The "Tiempos" sheet had two side-by-side lookup tables (5-step granularity in
# columns P:R and a 10-step granularity in columns T:V). The 5-step table
# (P:R) is unused/redundant, so select it and delete the entire columns,
# which shifts the remaining T:V table left to Q:S.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Tiempos")
$ws.Activate()
$ws.Range("P14").Select()

$ws.Range("P:R").Delete()
